# Insert a new weekly price record at row 50, shifting the existing
# rows 50-132 down to 51-133 (dimension grows from A1:R132 to A1:R133).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 50..132 down by one to make room for the new record.
$ws.Rows("50").Insert()

# Populate the newly inserted row 50 with the new data record.
$ws.Range("A50").Value = 4
$ws.Range("B50").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C50").Value = 'Los Lagos'
$ws.Range("D50").Value = 44757
$ws.Range("E50").Value = 10
$ws.Range("F50").Value = 100112052
$ws.Range("G50").Value = 'Albahaca'
$ws.Range("H50").Value = 'Sin especificar'
$ws.Range("I50").Value = 'Primera'
$ws.Range("J50").Value = 90
$ws.Range("K50").Value = 6500
$ws.Range("L50").Value = 6500
$ws.Range("M50").Value = 6500
$ws.Range("N50").Value = '$/paquete'
$ws.Range("O50").Value = 'Región de Arica y Parinacota'
$ws.Range("P50").Value = 6500
$ws.Range("Q50").Value = 1
$ws.Range("R50").Value = 'Hortaliza'
